$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = '61.983.30'
    $ws.Range("E2").Value = '  +1.61%  '
    $ws.Range("D3").Value = '2.410.78'
    $ws.Range("E3").Value = '  +1.63%  '
    $ws.Range("E4").Value = '  -0.06%  '
    $ws.Range("D5").Value = "'556.40"
    $ws.Range("D5").Style = "Normal"
    $ws.Range("D6").Value = "'142.40"
    $ws.Range("D6").Style = "Normal"
    $ws.Range("E6").Value = '  +3.41%  '
    $ws.Range("E7").Value = '  +0.10%  '
    $ws.Range("E8").Value = '  +0.71%  '
    $ws.Range("D9").Value = '2.405.17'
    $ws.Range("E9").Value = '  +1.37%  '
    $ws.Range("E10").Value = '  +1.25%  '
    $ws.Range("E11").Value = '  -0.95%  '
    $ws.Range("E12").Value = '  +1.24%  '
    $ws.Range("E13").Value = '  +1.40%  '
    $ws.Range("D14").Value = "'25.99"
    $ws.Range("D14").Style = "Normal"
    $ws.Range("E14").Value = '  +4.18%  '
    $ws.Range("E15").Value = '  +4.84%  '
    $ws.Range("D16").Value = '2.846.80'
    $ws.Range("E16").Value = '  +2.29%  '
    $ws.Range("D17").Value = '61.992.57'
    $ws.Range("E17").Value = '  +1.74%  '
    $ws.Range("D18").Value = '2.408.23'
    $ws.Range("E18").Value = '  +0.81%  '
    $ws.Range("E19").Value = '  +3.15%  '
    $ws.Range("E20").Value = '  +1.16%  '
    $ws.Range("D21").Value = "'322.82"
    $ws.Range("D21").Style = "Normal"
    $ws.Range("E21").Value = '  +0.99%  '
    $ws.Range("D22").Value = "'6.72"
    $ws.Range("D22").Style = "Normal"
    $ws.Range("E22").Value = '  +0.61%  '
    $ws.Range("E23").Value = '  +0.00%  '
    $ws.Range("D24").Value = "'65.28"
    $ws.Range("D24").Style = "Normal"
    $ws.Range("E24").Value = '  +1.73%  '
    $ws.Range("E25").Value = '  +2.12%  '
    $ws.Range("D26").Value = "'9.19"
    $ws.Range("D26").Style = "Normal"
    $ws.Range("E26").Value = '  +9.27%  '
    $ws.Range("D27").Value = "'578.72"
    $ws.Range("D27").Style = "Normal"
    $ws.Range("E27").Value = '  +15.02%  '
    $ws.Range("E28").Value = '  +0.15%  '
    $ws.Range("D29").Value = '2.528.00'
    $ws.Range("E29").Value = '  +1.71%  '
    $ws.Range("D30").Value = '0.0₃0928'
    $ws.Range("E30").Value = '  +5.83%  '
    $ws.Range("D32").Value = "'1.44"
    $ws.Range("D32").Style = "Normal"
    $ws.Range("E32").Value = '  +4.85%  '
    $ws.Range("E33").Value = '  -1.42%  '
    $ws.Range("E34").Value = '  +2.57%  '
    $ws.Range("E35").Value = '  +3.18%  '
    $ws.Range("E36").Value = '  +0.13%  '
    $ws.Range("D37").Value = "'5.65"
    $ws.Range("D37").Style = "Normal"
    $ws.Range("E37").Value = '  +6.07%  '
    $ws.Range("D38").Value = "'4.74"
    $ws.Range("D38").Style = "Normal"
    $ws.Range("E38").Value = '  +1.89%  '
    $ws.Range("E39").Value = '  +1.52%  '
    $ws.Range("B40").Value = 'EthereumClassic'
    $ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    $ws.Range("D40").Value = "'18.65"
    $ws.Range("D40").Style = "Normal"
    $ws.Range("E40").Value = '  +0.71%  '
    $ws.Range("B41").Value = 'Monero'
    $ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    $ws.Range("D41").Value = "'150.33"
    $ws.Range("D41").Style = "Normal"
    $ws.Range("E41").Value = '  +3.18%  '
    $ws.Range("D42").Value = "'1.82"
    $ws.Range("D42").Style = "Normal"
    $ws.Range("E42").Value = '  -2.56%  '
    $ws.Range("D43").Value = "'0.999"
    $ws.Range("D43").Style = "Normal"
    $ws.Range("D44").Value = "'2.30"
    $ws.Range("D44").Style = "Normal"
    $ws.Range("E44").Value = '  +12.53%  '
    $ws.Range("D45").Value = "'150.66"
    $ws.Range("D45").Style = "Normal"
    $ws.Range("E45").Value = '  +2.28%  '
    $ws.Range("E46").Value = '  +1.48%  '
    $ws.Range("D47").Value = "'0.0539"
    $ws.Range("D47").Style = "Normal"
    $ws.Range("E47").Value = '  +3.70%  '
    $ws.Range("D48").Value = "'20.15"
    $ws.Range("D48").Style = "Normal"
    $ws.Range("E48").Value = '  +4.97%  '
    $ws.Range("E49").Value = '  +2.38%  '
    $ws.Range("D50").Value = "'0.0921"
    $ws.Range("D50").Style = "Normal"
    $ws.Range("E50").Value = '  +1.33%  '
    $ws.Range("E51").Value = '  +2.29%  '
